$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update synthetics_from_this_seed (W), last_delta (X), last_neighbor_z (Y), timestamp (Z)
# for rows 2-60 per corrected pcsmote run results.

$ws.Cells.Item(2, 23).Value = 1
$ws.Cells.Item(2, 24).Value = 0.5361410903109534
$ws.Cells.Item(2, 25).Value = 80
$ws.Cells.Item(2, 26).Value = "2025-10-29T23:40:43.527204"

$ws.Cells.Item(3, 26).Value = "2025-10-29T23:40:43.527204"

$ws.Cells.Item(4, 26).Value = "2025-10-29T23:40:43.527204"

$ws.Cells.Item(5, 26).Value = "2025-10-29T23:40:43.527204"

$ws.Cells.Item(6, 23).Value = 2
$ws.Cells.Item(6, 24).Value = 0.5995480970097884
$ws.Cells.Item(6, 25).Value = 12
$ws.Cells.Item(6, 26).Value = "2025-10-29T23:40:43.527204"

$ws.Cells.Item(7, 23).Value = 1
$ws.Cells.Item(7, 24).Value = 0.4161706652665431
$ws.Cells.Item(7, 25).Value = 4
$ws.Cells.Item(7, 26).Value = "2025-10-29T23:40:43.527204"

$ws.Cells.Item(8, 23).Value = 3
$ws.Cells.Item(8, 24).Value = 0.5424541179848884
$ws.Cells.Item(8, 25).Value = 107
$ws.Cells.Item(8, 26).Value = "2025-10-29T23:40:43.528204"

$ws.Cells.Item(9, 23).Value = 5
$ws.Cells.Item(9, 24).Value = 0.4834822006297558
$ws.Cells.Item(9, 26).Value = "2025-10-29T23:40:43.528204"

$ws.Cells.Item(10, 23).Value = 3
$ws.Cells.Item(10, 24).Value = 0.4013904261062382
$ws.Cells.Item(10, 25).Value = 50
$ws.Cells.Item(10, 26).Value = "2025-10-29T23:40:43.531685"

$ws.Cells.Item(11, 26).Value = "2025-10-29T23:40:43.531685"

$ws.Cells.Item(12, 23).Value = 3
$ws.Cells.Item(12, 24).Value = 0.4244175909401347
$ws.Cells.Item(12, 25).Value = 120
$ws.Cells.Item(12, 26).Value = "2025-10-29T23:40:43.531685"

$ws.Cells.Item(13, 23).Value = 3
$ws.Cells.Item(13, 24).Value = 0.4646405864041511
$ws.Cells.Item(13, 25).Value = 65
$ws.Cells.Item(13, 26).Value = "2025-10-29T23:40:43.531685"

$ws.Cells.Item(14, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(15, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(16, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(17, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(18, 23).Value = 1
$ws.Cells.Item(18, 24).Value = 0.5924894589884222
$ws.Cells.Item(18, 25).Value = 140
$ws.Cells.Item(18, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(19, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(20, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(21, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(22, 23).Value = 1
$ws.Cells.Item(22, 24).Value = 0.4895566329146183
$ws.Cells.Item(22, 25).Value = 58
$ws.Cells.Item(22, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(23, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(24, 26).Value = "2025-10-29T23:40:43.532681"

$ws.Cells.Item(25, 23).Value = 3
$ws.Cells.Item(25, 24).Value = 0.4807672342116082
$ws.Cells.Item(25, 25).Value = 6
$ws.Cells.Item(25, 26).Value = "2025-10-29T23:40:43.533679"

$ws.Cells.Item(26, 23).Value = 3
$ws.Cells.Item(26, 24).Value = 0.4544264498769271
$ws.Cells.Item(26, 25).Value = 146
$ws.Cells.Item(26, 26).Value = "2025-10-29T23:40:43.533679"

$ws.Cells.Item(27, 26).Value = "2025-10-29T23:40:43.533679"

$ws.Cells.Item(28, 23).Value = 4
$ws.Cells.Item(28, 24).Value = 0.5953229911665305
$ws.Cells.Item(28, 25).Value = 134
$ws.Cells.Item(28, 26).Value = "2025-10-29T23:40:43.533679"

$ws.Cells.Item(29, 23).Value = 5
$ws.Cells.Item(29, 24).Value = 0.4690142496053366
$ws.Cells.Item(29, 25).Value = 102
$ws.Cells.Item(29, 26).Value = "2025-10-29T23:40:43.533679"

$ws.Cells.Item(30, 23).Value = 9
$ws.Cells.Item(30, 24).Value = 0.5797108377054159
$ws.Cells.Item(30, 25).Value = 23
$ws.Cells.Item(30, 26).Value = "2025-10-29T23:40:43.556941"

$ws.Cells.Item(31, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(32, 23).Value = 9
$ws.Cells.Item(32, 24).Value = 0.5712648583756185
$ws.Cells.Item(32, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(33, 23).Value = 5
$ws.Cells.Item(33, 24).Value = 0.5381875476204931
$ws.Cells.Item(33, 25).Value = 6
$ws.Cells.Item(33, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(34, 23).Value = 6
$ws.Cells.Item(34, 24).Value = 0.576527268637868
$ws.Cells.Item(34, 25).Value = 112
$ws.Cells.Item(34, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(35, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(36, 23).Value = 7
$ws.Cells.Item(36, 24).Value = 0.5284063292308575
$ws.Cells.Item(36, 25).Value = 61
$ws.Cells.Item(36, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(37, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(38, 23).Value = 9
$ws.Cells.Item(38, 24).Value = 0.4557742705184364
$ws.Cells.Item(38, 25).Value = 11
$ws.Cells.Item(38, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(39, 26).Value = "2025-10-29T23:40:43.557941"

$ws.Cells.Item(40, 23).Value = 14
$ws.Cells.Item(40, 24).Value = 0.4484319876554852
$ws.Cells.Item(40, 25).Value = 99
$ws.Cells.Item(40, 26).Value = "2025-10-29T23:40:43.558987"

$ws.Cells.Item(41, 23).Value = 7
$ws.Cells.Item(41, 24).Value = 0.4418143241475427
$ws.Cells.Item(41, 25).Value = 134
$ws.Cells.Item(41, 26).Value = "2025-10-29T23:40:43.558987"

$ws.Cells.Item(42, 23).Value = 7
$ws.Cells.Item(42, 24).Value = 0.5033271782542028
$ws.Cells.Item(42, 25).Value = 3
$ws.Cells.Item(42, 26).Value = "2025-10-29T23:40:43.559518"

$ws.Cells.Item(43, 23).Value = 2
$ws.Cells.Item(43, 24).Value = 0.4641560129943472
$ws.Cells.Item(43, 25).Value = 50
$ws.Cells.Item(43, 26).Value = "2025-10-29T23:40:43.594728"

$ws.Cells.Item(44, 26).Value = "2025-10-29T23:40:43.594728"

$ws.Cells.Item(45, 26).Value = "2025-10-29T23:40:43.595730"

$ws.Cells.Item(46, 23).Value = 5
$ws.Cells.Item(46, 24).Value = 0.5881046528979208
$ws.Cells.Item(46, 25).Value = 146
$ws.Cells.Item(46, 26).Value = "2025-10-29T23:40:43.595730"

$ws.Cells.Item(47, 23).Value = 3
$ws.Cells.Item(47, 24).Value = 0.4895566329146183
$ws.Cells.Item(47, 25).Value = 23
$ws.Cells.Item(47, 26).Value = "2025-10-29T23:40:43.595730"

$ws.Cells.Item(48, 23).Value = 4
$ws.Cells.Item(48, 24).Value = 0.571671760962744
$ws.Cells.Item(48, 25).Value = 61
$ws.Cells.Item(48, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(49, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(50, 23).Value = 9
$ws.Cells.Item(50, 24).Value = 0.5675420211814656
$ws.Cells.Item(50, 25).Value = 118
$ws.Cells.Item(50, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(51, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(52, 23).Value = 4
$ws.Cells.Item(52, 24).Value = 0.561900209227943
$ws.Cells.Item(52, 25).Value = 121
$ws.Cells.Item(52, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(53, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(54, 23).Value = 6
$ws.Cells.Item(54, 24).Value = 0.5953229911665305
$ws.Cells.Item(54, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(55, 26).Value = "2025-10-29T23:40:43.596732"

$ws.Cells.Item(56, 26).Value = "2025-10-29T23:40:43.597731"

$ws.Cells.Item(57, 23).Value = 10
$ws.Cells.Item(57, 24).Value = 0.4440482095131097
$ws.Cells.Item(57, 25).Value = 153
$ws.Cells.Item(57, 26).Value = "2025-10-29T23:40:43.597731"

$ws.Cells.Item(58, 26).Value = "2025-10-29T23:40:43.597731"

$ws.Cells.Item(59, 23).Value = 5
$ws.Cells.Item(59, 24).Value = 0.450783082786869
$ws.Cells.Item(59, 26).Value = "2025-10-29T23:40:43.597731"

$ws.Cells.Item(60, 26).Value = "2025-10-29T23:40:43.597731"

# Remove obsolete rows 61-68 (window for pureza_proporcion sample increased,
# dataset now only spans rows 2-60) and fix up conditional formatting range.
$ws.Range("A61:Z68").Delete()

$cf = $ws.Range("A2:Z68").FormatConditions
if ($cf.Count -ge 1) {
    $cf.Item(1).ModifyAppliesToRange($ws.Range("A2:Z60"))
}
